$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text/percent/name/url assignments (unambiguous, never auto-converted to numbers)
$ws.Range("D2").Value = '65.033.66'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '3.525.63'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").Value = '3.523.86'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("E11").Value = '  +3.27%  '
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '4.122.12'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '3.523.06'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '64.965.54'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("E21").Value = '  -1.77%  '
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.660.29'
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -3.22%  '
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("E29").Value = '  +9.52%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("D33").Value = '3.517.80'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("E37").Value = '  +6.77%  '
$ws.Range("E38").Value = '  +3.37%  '
$ws.Range("E39").Value = '  +1.07%  '
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  +5.77%  '
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("D50").Value = '2.415.30'
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("E51").Value = '  +6.17%  '

# Price column values that look like plain numbers need to be forced to remain text,
# otherwise Excel auto-converts them to numeric values (losing formatting/precision).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.49'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.29'
$ws.Range("D6").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.15'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.69'
$ws.Range("D14").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.16'
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.47'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.72'
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '393.42'
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.580'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.81'
$ws.Range("D24").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.37'
$ws.Range("D32").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.14'
$ws.Range("D34").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.00'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '168.57'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0816'
$ws.Range("D41").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.97'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.94'
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.45'
$ws.Range("D47").ClearFormats()
